# "uploaded just snap caps" - add the new snap-cap color variants to the
# products_URLimages worksheet, each with a name (cols A & B) and an
# image URL (col C) that is turned into a hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A got narrower in this revision.
$ws.Columns.Item(1).ColumnWidth = 44.5

# The existing last row (152) gets a real hyperlink + the "Hyperlink" cell
# style applied to its URL cell, matching the rest of the sheet.
$lastRow = 152
$lastUrl = $ws.Cells.Item($lastRow, 3).Value()
$ws.Hyperlinks.Add($ws.Cells.Item($lastRow, 3), $lastUrl) | Out-Null
$ws.Cells.Item($lastRow, 3).Style = "Hyperlink"

# New snap-cap products to append below the existing data.
$items = @(
    @("snap-cap-black",  "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-black.jpg"),
    @("snap-cap-blue",   "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-blue.jpg"),
    @("snap-cap-gray",   "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-gray.jpg"),
    @("snap-cap-green",  "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-green.jpg"),
    @("snap-cap-pink",   "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-pink.jpg"),
    @("snap-cap-red",    "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-red.jpg"),
    @("snap-cap-yellow", "https://raw.githubusercontent.com/htslabs/images/main/snap-cap-yellow.jpg")
)

$row = $lastRow + 1
foreach ($item in $items) {
    $name = $item[0]
    $url = $item[1]

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $url

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 3), $url) | Out-Null
    $ws.Cells.Item($row, 3).Style = "Hyperlink"

    $row = $row + 1
}

# Reflect the new selection / scrolled position left behind by the edit.
$ws.Range("A153:XFD159").Select() | Out-Null
